# download_request_example.xlsx — add un-archiving / copy-destination columns
# and turn the source URL into a real hyperlink.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename/retext a couple of cells (order matters for shared-string table order).
$ws.Range("D2").Value = "{box}/methylation"
$ws.Range("D1").Value = "destination_path"
$ws.Range("E2").Value = "Yes"

# Turn the GoogleDrive URL in B2 into a clickable hyperlink (keeps display text).
$url = $ws.Range("B2").Value()
$ws.Hyperlinks.Add($ws.Range("B2"), $url) | Out-Null

# Move the selection down to the (empty) row below the data, selecting the whole row.
$ws.Rows(3).Select()
